$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated values in row 2 (Beta)
$ws.Range("F2").Value = 9.558705292934931
$ws.Range("G2").Value = 7.719402956952014
$ws.Range("H2").Value = 11.52883248399174
$ws.Range("I2").Value = 1.679789493733906
$ws.Range("J2").Value = 0.8597426907656701
$ws.Range("K2").Value = 2.650064815982796
$ws.Range("L2").Value = 0.142662760886534
$ws.Range("M2").Value = 0.0889265508633756
$ws.Range("N2").Value = 0.2110483859013716

# Update recalculated values in row 3 (Gamma)
$ws.Range("F3").Value = 0.3716748077026123
$ws.Range("G3").Value = 0.0001679632586709231
$ws.Range("H3").Value = 1.064023942517129
$ws.Range("I3").Value = 0.3048827760661708
$ws.Range("J3").Value = 0.0001432833411799475
$ws.Range("K3").Value = 0.8721313737180489
$ws.Range("L3").Value = 0.373158765322302
$ws.Range("M3").Value = 0.0001698325457199391
$ws.Range("N3").Value = 1.068595975283141

# Add new row 4 (Beta + Gamma)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 9.930380100637544
$ws.Range("G4").Value = 7.719570920210685
$ws.Range("H4").Value = 12.59285642650887
$ws.Range("I4").Value = 1.984672269800077
$ws.Range("J4").Value = 0.85988597410685
$ws.Range("K4").Value = 3.522196189700845
$ws.Range("L4").Value = 0.5158215262088359
$ws.Range("M4").Value = 0.08909638340909554
$ws.Range("N4").Value = 1.279644361184513
